$wb = $excel.ActiveWorkbook

# --- Sheet2: Withdraw History ---
$ws2 = $wb.Worksheets.Item("Withdraw History")
$ws2.Range("A1").Value = "Amount"
$ws2.Range("B1").Value = "Time"
$ws2.Range("C1").Value = "Day"
$ws2.Range("D1").Value = "Month"
$ws2.Range("E1").Value = "Year"
$ws2.Range("F1").Value = "Place"
$ws2.Range("A1:F1").Select() | Out-Null

# --- Sheet3: Deposit History ---
$ws3 = $wb.Worksheets.Item("Deposit History")
$ws3.Range("A1").Value = "Amount"
$ws3.Range("B1").Value = "Time"
$ws3.Range("C1").Value = "Day"
$ws3.Range("D1").Value = "Month"
$ws3.Range("E1").Value = "Year"
$ws3.Range("F1").Value = "Place"
$ws3.Range("A1:F1").Select() | Out-Null

# --- Sheet4: Transfer History ---
$ws4 = $wb.Worksheets.Item("Transfer History")
$ws4.Range("A1").Value = "Amount"
$ws4.Range("B1").Value = "Time"
$ws4.Range("C1").Value = "Day"
$ws4.Range("D1").Value = "Month"
$ws4.Range("E1").Value = "Year"
$ws4.Range("F1").Value = "Person"
$ws4.Range("A1:E1").Select() | Out-Null

# --- Sheet5: Absolute History ---
$ws5 = $wb.Worksheets.Item("Absolute History")
$ws5.Range("A1").Value = "Amount"
$ws5.Range("B1").Value = "Time"
$ws5.Range("C1").Value = "Day"
$ws5.Range("D1").Value = "Month"
$ws5.Range("E1").Value = "Year"
$ws5.Range("F1").Value = "Place/Person"
$ws5.Range("A2").Select() | Out-Null

# --- Sheet1: Amount becomes the active tab/sheet ---
$ws1 = $wb.Worksheets.Item("Amount")
$ws1.Activate()

Write-Host "Applied skeleton headers"
